# Word COM-interop script implementing the "Updated contents with some
# more reasoning." commit.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Scaling bullet list: append extra reasoning to each of the three
#    bullet points about multi-instance / multi-container / Kubernetes
#    deployments.
# ---------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("- deploying multiple instances of Tornado exploiting multiple threads,", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" but yet it would requires a shared memory approach (e.g. using memcache).")

$rng = $d.Content
$rng.Find.Execute("- deploying multiple containers,", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" and using a load balancer in front of them, but yet a DB shared-access logic should be implemented (e.g. to avoid possible double mongodb " + [char]8220 + "_id" + [char]8221 + " a shared list can be used by all the instances).")

$rng = $d.Content
$rng.Find.Execute("- deploying Tornado on Kubernetes (e.g. using pod replicas)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" and a load balancer in front of them.")

# ---------------------------------------------------------------------
# 2) Monitoring paragraph: the three runs that make up this sentence are
#    simply merged back into a single contiguous sentence (no content
#    change), so a self-replace collapses them.
# ---------------------------------------------------------------------

$monitoringText = "Monitoring a Kubernetes cluster could be the easiest way to go (e.g. deploying prometheus with helm chart using a daemonset)."
$rng = $d.Content
$rng.Find.Execute($monitoringText, `
    $true, $false, $false, $false, $false, $true, 1, $false, $monitoringText, 2)

# ---------------------------------------------------------------------
# 3) High availability paragraph: likewise, the "It is coupled with ...
#    VMs/Nodes locations ..." sentence is reassembled from several runs
#    into one contiguous run.
# ---------------------------------------------------------------------

$haText = "It is coupled with the scaling approach, caring about the VMs/Nodes locations " + [char]8211 + " interesting for an international deployments, less interesting for a national scale one (since there are not many Azure datacenters in Italy, unless we want to proceed with an hybrid cloud/on-prem approach)."
$rng = $d.Content
$rng.Find.Execute($haText, `
    $true, $false, $false, $false, $false, $true, 1, $false, $haText, 2)

# ---------------------------------------------------------------------
# 4) Jmeter paragraph: "collection" is struck through and replaced by
#    "document", and a comment is attached explaining the correction.
# ---------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("collection", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.StrikeThrough = 1

$rng = $d.Content
$rng.Find.Execute("collection ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter("document ")

$word.UserName = "Unknown Author"

$rng = $d.Content
$rng.Find.Execute("the MongoDB collection document max.16mb size", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$commentText = "This test doesn" + [char]8217 + "t make much sense. I mistakenly thought that " + [char]8220 + "uncapped" + [char]8221 + " collection may have a fixed size. The fixed size is only for documents (namely JSON objects)."
$comment = $d.Comments.Add($rng, $commentText)
$comment.Initial = ""
